# Apply the FlashScore odds updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("AE2").Value = 29
$ws.Range("AX2").Value = 41

# Row 3
$ws.Range("G3").Value = 1.45
$ws.Range("I3").Value = 7.5
$ws.Range("J3").Value = 2
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("AW3").Value = 8
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 351

# Row 5
$ws.Range("G5").Value = 2.7
$ws.Range("I5").Value = 2.88
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 1.91
$ws.Range("L5").Value = 3.6
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("W5").Value = 7
$ws.Range("X5").Value = 12
$ws.Range("AA5").Value = 26
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 12
$ws.Range("AO5").Value = 17
$ws.Range("AS5").Value = 301
$ws.Range("AY5").Value = 29

$wb.Save()
